$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Renumber the starting page for this section: 161 -> 129 ---------------
# This rewrites <w:pgNumType w:start="161"/> to <w:pgNumType w:start="129"/>
# in the section properties.
$defaultFooter = $sec.Footers(1)
$defaultFooter.PageNumbers.StartingNumber = 129

# --- Fix up the cached PAGE field result shown in the default footer -------
# The footer contains a " PAGE \* MERGEFORMAT " field whose last computed
# result ("161") is cached as literal text in the document; update that
# cached text to match the new starting page number ("129").
$footerRange = $defaultFooter.Range
$chars = $footerRange.Characters
$cachedPageNumber = $footerRange.Duplicate
$cachedPageNumber.Start = $chars.Item(1).Start
$cachedPageNumber.End = $chars.Item(3).End
if ($cachedPageNumber.Text -eq "161") {
    $cachedPageNumber.Text = "129"
}

# --- Drop the even/first-page headers & footers -----------------------------
# The section only needs a single default footer (no headers, no even/first
# page header or footer). Clear their content and make sure the "different
# odd/even" and "different first page" behaviour is switched off so the
# section no longer distinguishes those variants.
for ($i = 2; $i -le 3; $i++) {
    $sec.Headers($i).Range.Text = ""
    $sec.Footers($i).Range.Text = ""
}
$sec.Headers(1).Range.Text = ""

$sec.PageSetup.OddAndEvenPagesHeaderFooter = $true
$sec.PageSetup.OddAndEvenPagesHeaderFooter = $false
$sec.PageSetup.DifferentFirstPageHeaderFooter = $true
$sec.PageSetup.DifferentFirstPageHeaderFooter = $false

Write-Output "done"
